{"js": "// Load all paragraphs in the document body so we can locate the\n// \"3.1 ...\" / \"3.2 ...\" task items by their text.\nconst paragraphs = context.document.body.paragraphs;\nparagraphs.load(\"items/text\");\nawait context.sync();\n\nlet para31 = null;\nlet para32 = null;\n\nfor (let i = 0; i < paragraphs.items.length; i++) {\n  const text = paragraphs.items[i].text;\n  if (text.indexOf(\"3.1 \u0441\u0434\u0435\u043b\u0430\u0442\u044c \u043a\u043d\u043e\u043f\u043a\u0443 \u00ab\u0434\u043e\u0431\u0430\u0432\u0438\u0442\u044c \u0432 \u043a\u043e\u0440\u0437\u0438\u043d\u0443\u00bb\") === 0) {\n    para31 = paragraphs.items[i];\n  } else if (text.indexOf(\"3.2 \u0441\u043e\u0437\u0434\u0430\u0442\u044c \u043c\u0430\u0441c\u0438\u0432 \u0434\u043b\u044f \u043a\u043e\u0440\u0437\u0438\u043d\u044b \u0447\u0435\u0440\u0435\u0437 localStorage\") === 0) {\n    para32 = paragraphs.items[i];\n  }\n}\n\nif (!para31 || !para32) {\n  throw new Error(\"Could not locate target paragraphs 3.1 / 3.2\");\n}\n\n// 3.1 -> append \" //\" to the end of the line.\npara31.insertText(\" //\", \"End\");\n\n// 3.2 -> append \"//\" (no leading space) to the end of the line.\npara32.insertText(\"//\", \"End\");\n\n// Insert the two new task paragraphs right after \"3.2 ...\", preserving the\n// same paragraph/run formatting (inherited automatically by insertParagraph).\nconst para33 = para32.insertParagraph(\"3.3 \u0441\u043e\u0437\u0434\u0430\u0442\u044c \u0432\u0441\u043f\u043b\u044b\u0432\u0430\u044e\u0449\u0435\u0435 \u043e\u043a\u043d\u043e \u0434\u043b\u044f \u043a\u043e\u0440\u0437\u0438\u043d\u044b\", \"After\");\npara33.insertParagraph(\"3.3.1 \u0434\u043e\u043f\u0438\u0441\u0430\u0442\u044c \u0432\u0451\u0440\u0441\u0442\u043a\u0443 \u0434\u043b\u044f \u0441\u0443\u0448\u0438 \u0432 \u043a\u043e\u0440\u0437\u0438\u043d\u0435\", \"After\");\n\nawait context.sync();\n", "ps1": "$d = $word.ActiveDocument\n\n# 1) Append \" //\" to the end of the \"3.1 ...\" task line.\n$rng1 = $d.Content\n$rng1.Find.Execute(\"3.1 \u0441\u0434\u0435\u043b\u0430\u0442\u044c \u043a\u043d\u043e\u043f\u043a\u0443 \u00ab\u0434\u043e\u0431\u0430\u0432\u0438\u0442\u044c \u0432 \u043a\u043e\u0440\u0437\u0438\u043d\u0443\u00bb\") | Out-Null\n$rng1.Collapse(0)\n$rng1.InsertAfter(\" //\")\n\n# 2) Append \"//\" (no leading space) to the end of the \"3.2 ...\" task line.\n$rng2 = $d.Content\n$rng2.Find.Execute(\"3.2 \u0441\u043e\u0437\u0434\u0430\u0442\u044c \u043c\u0430\u0441c\u0438\u0432 \u0434\u043b\u044f \u043a\u043e\u0440\u0437\u0438\u043d\u044b \u0447\u0435\u0440\u0435\u0437 localStorage\") | Out-Null\n$rng2.Collapse(0)\n$rng2.InsertAfter(\"//\")\n\n# 3) Insert a new \"3.3 ...\" paragraph right after the (now updated) \"3.2 ...\" line,\n#    inheriting its paragraph/run formatting.\n$rng3 = $d.Content\n$rng3.Find.Execute(\"3.2 \u0441\u043e\u0437\u0434\u0430\u0442\u044c \u043c\u0430\u0441c\u0438\u0432 \u0434\u043b\u044f \u043a\u043e\u0440\u0437\u0438\u043d\u044b \u0447\u0435\u0440\u0435\u0437 localStorage//\") | Out-Null\n$para32 = $rng3.Paragraphs(1)\n$para32.Range.InsertParagraphAfter()\n$para33 = $para32.Next()\n$para33.Range.InsertAfter(\"3.3 \u0441\u043e\u0437\u0434\u0430\u0442\u044c \u0432\u0441\u043f\u043b\u044b\u0432\u0430\u044e\u0449\u0435\u0435 \u043e\u043a\u043d\u043e \u0434\u043b\u044f \u043a\u043e\u0440\u0437\u0438\u043d\u044b\")\n\n# 4) Insert a new \"3.3.1 ...\" paragraph right after the \"3.3 ...\" line.\n$para33.Range.InsertParagraphAfter()\n$para331 = $para33.Next()\n$para331.Range.InsertAfter(\"3.3.1 \u0434\u043e\u043f\u0438\u0441\u0430\u0442\u044c \u0432\u0451\u0440\u0441\u0442\u043a\u0443 \u0434\u043b\u044f \u0441\u0443\u0448\u0438 \u0432 \u043a\u043e\u0440\u0437\u0438\u043d\u0435\")\n"}
